$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1954.3636
$ws.Range("J41").Value = 99.5
$ws.Range("L41").Value = 99.5
$ws.Range("N41").Value = -979.5

$ws.Range("H98").Value = 918.6667
$ws.Range("I98").Value = 918.6667
$ws.Range("K98").Value = 918.6667
$ws.Range("M98").Value = 579.3333

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("N100").Value = 0
$ws.Range("L100").ClearContents()
$ws.Range("M100").ClearContents()

$ws.Range("H116").Value = 3664.6667
$ws.Range("I116").Value = 3664.6667
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 3664.6667
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = -222.6667000000002
$ws.Range("M116").ClearContents()

$ws.Range("H122").Value = 918.6667
$ws.Range("I122").Value = 918.6667
$ws.Range("K122").Value = 2756.0001
$ws.Range("M122").Value = -306.0001000000002

$ws.Range("H123").Value = 30449.5
$ws.Range("J123").Value = 30449.5
$ws.Range("L123").Value = 30449.5
$ws.Range("N123").Value = -40249.5

$ws.Range("H132").Value = 2700
$ws.Range("I132").Value = 2700
$ws.Range("K132").Value = 8100
$ws.Range("M132").Value = -5570

$ws.Range("H135").Value = 1614.5385
$ws.Range("I135").Value = 1729.1
$ws.Range("J135").Value = 1232.6666
$ws.Range("K135").Value = 15561.9
$ws.Range("L135").Value = 11093.9994
$ws.Range("M135").Value = -13026.9
$ws.Range("N135").Value = -16163.9994

$ws.Range("H138").Value = 3164.8057
$ws.Range("I138").Value = 1547.0714
$ws.Range("J138").Value = 4194.273
$ws.Range("K138").Value = 4641.2142
$ws.Range("L138").Value = 12582.819
$ws.Range("M138").Value = 498.7857999999997
$ws.Range("N138").Value = -22862.819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1265.2222
$ws.Range("I5").Value = 1462.5714
$ws.Range("J5").Value = 574.5
$ws.Range("K5").Value = 1462.5714
$ws.Range("L5").Value = 574.5
$ws.Range("M5").Value = -1350.5714
$ws.Range("N5").Value = -798.5

$ws.Range("H102").Value = 2149.5
$ws.Range("I102").Value = 2049
$ws.Range("J102").Value = 2250
$ws.Range("K102").Value = 2049
$ws.Range("L102").Value = 2250
$ws.Range("M102").Value = -427
$ws.Range("N102").Value = -5494

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1265.2222
$ws.Range("I4").Value = 1462.5714
$ws.Range("J4").Value = 574.5
$ws.Range("K4").Value = 1462.5714
$ws.Range("L4").Value = 574.5
$ws.Range("M4").Value = -1347.5714
$ws.Range("N4").Value = -804.5

$ws.Range("H20").Value = 4425.5
$ws.Range("I20").Value = 10008
$ws.Range("J20").Value = 2564.6667
$ws.Range("K20").Value = 10008
$ws.Range("L20").Value = 2564.6667
$ws.Range("M20").Value = -9761
$ws.Range("N20").Value = -3058.6667

$ws.Range("H105").Value = 2267.6
$ws.Range("I105").Value = 2267.6
$ws.Range("K105").Value = 2267.6
$ws.Range("M105").Value = -520.5999999999999

$ws.Range("H107").Value = 3077.1428
$ws.Range("I107").Value = 3210
$ws.Range("K107").Value = 3210
$ws.Range("M107").Value = -1290

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 189474.83
$ws.Range("J94").Value = 7540
$ws.Range("L94").Value = 7540
$ws.Range("N94").Value = -8442

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("N99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").ClearContents()

$ws.Range("H105").Value = 3607.4285
$ws.Range("I105").Value = 3063.125
$ws.Range("J105").Value = 4333.1665
$ws.Range("K105").Value = 3063.125
$ws.Range("L105").Value = 4333.1665
$ws.Range("M105").Value = -1316.125
$ws.Range("N105").Value = -7827.1665

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("N122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("N126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 4614.4165
$ws.Range("I132").Value = 4820.2
$ws.Range("K132").Value = 14460.6
$ws.Range("M132").Value = -11930.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6914.5713
$ws.Range("I70").Value = 5599.25
$ws.Range("J70").Value = 8668.333000000001
$ws.Range("K70").Value = 5599.25
$ws.Range("L70").Value = 8668.333000000001
$ws.Range("M70").Value = -5329.25
$ws.Range("N70").Value = -9208.333000000001

$ws.Range("H73").Value = 6914.5713
$ws.Range("I73").Value = 5599.25
$ws.Range("J73").Value = 8668.333000000001
$ws.Range("K73").Value = 5599.25
$ws.Range("L73").Value = 8668.333000000001
$ws.Range("M73").Value = -4663.25
$ws.Range("N73").Value = -10540.333

$ws.Range("H113").Value = 3673.2
$ws.Range("I113").Value = 1471
$ws.Range("J113").Value = 5141.3335
$ws.Range("K113").Value = 1471
$ws.Range("L113").Value = 5141.3335
$ws.Range("M113").Value = 699
$ws.Range("N113").Value = -9481.333500000001

$ws.Range("H123").Value = 34835.715
$ws.Range("J123").Value = 34835.715
$ws.Range("L123").Value = 34835.715
$ws.Range("N123").Value = -39735.715

$ws.Range("H126").Value = 4430.3335
$ws.Range("I126").Value = 4398
$ws.Range("J126").Value = 4446.5
$ws.Range("K126").Value = 13194
$ws.Range("L126").Value = 13339.5
$ws.Range("M126").Value = -10724
$ws.Range("N126").Value = -18279.5

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("N128").Value = 0
$ws.Range("L128").ClearContents()

$ws.Range("H132").Value = 4658
$ws.Range("J132").Value = 5443.5
$ws.Range("L132").Value = 16330.5
$ws.Range("N132").Value = -21390.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2000.6666
$ws.Range("I61").Value = 2004
$ws.Range("K61").Value = 2004
$ws.Range("M61").Value = -1802

$ws.Range("H75").Value = 500173
$ws.Range("J75").Value = 500173
$ws.Range("L75").Value = 500173
$ws.Range("N75").Value = -502045

$ws.Range("H78").Value = 500173
$ws.Range("J78").Value = 500173
$ws.Range("L78").Value = 1500519
$ws.Range("N78").Value = -1509879

$ws.Range("H113").Value = 2000.6666
$ws.Range("I113").Value = 2004
$ws.Range("K113").Value = 2004
$ws.Range("M113").Value = 166

$ws.Range("H132").Value = 4293.8887
$ws.Range("I132").Value = 4666.3335
$ws.Range("J132").Value = 4107.6665
$ws.Range("K132").Value = 13999.0005
$ws.Range("L132").Value = 12322.9995
$ws.Range("M132").Value = -11469.0005
$ws.Range("N132").Value = -17382.9995

$ws.Range("H136").Value = 1707.9678
$ws.Range("I136").Value = 1523.6316
$ws.Range("K136").Value = 4570.8948
$ws.Range("M136").Value = -2020.8948

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1306.2
$ws.Range("I132").Value = 1306.2
$ws.Range("K132").Value = 3918.6
$ws.Range("M132").Value = -1388.6
